$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 6: Execute flag changes from Yes to No
$ws.Range("C6").Value = "No"

# Add new row 10: amazon test case
$ws.Range("A10").Value = "searchProduct"
$ws.Range("B10").Value = "validate the Sorting works"
$ws.Range("C10").Value = "Yes"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Extend the data validation list on column C to include the new row
$ws.Range("C2:C10").Validation.Delete()
$ws.Range("C2:C10").Validation.Add(3, 1, 1, """Yes,No""")
$ws.Range("C2:C10").Validation.IgnoreBlank = $true
$ws.Range("C2:C10").Validation.InCellDropdown = $true
$ws.Range("C2:C10").Validation.ShowInput = $true
$ws.Range("C2:C10").Validation.ShowError = $true

# Update selection to reflect the newly active cell
$ws.Range("A10").Select()
